$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format ("@") for Price (D) column cells so that numeric-looking
# strings (e.g. "335.80", "1.000") are preserved exactly as text, matching the
# original inline-string cell contents instead of being coerced into numbers.
$priceCells = @("D2","D3","D5","D7","D8","D9","D11","D13","D14","D15","D16","D18","D19","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped from coinranking.com
$ws.Range("D2").Value = "27.264.98"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "1.786.86"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "335.80"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "0.3812"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D8").Value = "0.3421"
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("D9").Value = "48.45"
$ws.Range("E9").Value = "  -3.20%  "
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").Value = "0.07488"
$ws.Range("E11").Value = "  -3.45%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "21.93"
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("D14").Value = "6.462"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").Value = "1.787.02"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "7.076"
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").Value = "0.06651"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").Value = "83.85"
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "6.625"
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("D22").Value = "17.35"
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").Value = "27.266.78"
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("D24").Value = "12.38"
$ws.Range("E24").Value = "  -6.26%  "
$ws.Range("D25").Value = "2.411"
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "2.542"
$ws.Range("E26").Value = "  -5.21%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "1.486"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").Value = "21.31"
$ws.Range("E28").Value = "  -3.58%  "
$ws.Range("D29").Value = "153.58"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "1.988.88"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("D31").Value = "134.16"
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").Value = "4.010"
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("D33").Value = "6.095"
$ws.Range("E33").Value = "  -4.56%  "
$ws.Range("D34").Value = "0.08712"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("D35").Value = "13.34"
$ws.Range("E35").Value = "  -4.32%  "
$ws.Range("D36").Value = "1.657"
$ws.Range("E36").Value = "  -3.63%  "
$ws.Range("D37").Value = "0.6943"
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("D38").Value = "5.451"
$ws.Range("E38").Value = "  -3.82%  "
$ws.Range("D39").Value = "0.2209"
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("D40").Value = "0.06338"
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("D41").Value = "8.809"
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("D42").Value = "0.02344"
$ws.Range("D43").Value = "1.238"
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("D44").Value = "14.40"
$ws.Range("E44").Value = "  -4.30%  "
$ws.Range("D45").Value = "0.6519"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "3.855"
$ws.Range("E47").Value = "  -4.80%  "
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("D49").Value = "129.12"
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("D50").Value = "0.07138"
$ws.Range("E50").Value = "  -3.05%  "
$ws.Range("D51").Value = "78.96"
$ws.Range("E51").Value = "  -2.33%  "
